$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# For data rows 2 through 23, clear out the columns AB..AK (10 cols) and AM,
# while leaving AL (PREVIOUS ACCOMPLISHMENT) untouched.
for ($r = 2; $r -le 23; $r++) {
    $ws.Range("AB$r`:AK$r").ClearContents()
    $ws.Range("AM$r").ClearContents()
}
